# MARS_m3c2_stats_clouds.xlsx - refresh the CloudStats header row to match the
# updated statistics-export pipeline (see commit message):
#  - "Mean NN Dist (1..6)" / "Mean Dist to 6-NN" relabelled to
#    "Mean NN Dist All" / "Mean NN Dist k-th"
#  - Local Density / Roughness headers drop their unit suffixes
#  - New geometric-feature columns are inserted (Anisotropy, Omnivariance,
#    Eigenentropy, Curvature, Verticality) ahead of the trailing
#    "Normal Std Angle [deg]" / "Radius [m]" / "k-NN" / "Sampled Points" /
#    "File" / "Folder" columns, which shift to the right to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final, full set of column headers for row 1, in column order (A, B, C, ...).
$headers = @(
    "Timestamp",
    "Num Points",
    "Area Source",
    "Area XY [m^2]",
    "Density Global [pt/m^2]",
    "Z Min",
    "Z Max",
    "Z Mean",
    "Z Median",
    "Z Std",
    "Z Q05",
    "Z Q25",
    "Z Q75",
    "Z Q95",
    "Mean NN Dist All",
    "Mean NN Dist k-th",
    "Local Density Mean",
    "Local Density Median",
    "Local Density Q05",
    "Local Density Q95",
    "Roughness Mean",
    "Roughness Median",
    "Roughness Q05",
    "Roughness Q95",
    "Linearity Mean",
    "Linearity Median",
    "Planarity Mean",
    "Planarity Median",
    "Sphericity Mean",
    "Sphericity Median",
    "Anisotropy Mean",
    "Anisotropy Median",
    "Omnivariance Mean",
    "Omnivariance Median",
    "Eigenentropy Mean",
    "Eigenentropy Median",
    "Curvature Mean",
    "Curvature Median",
    "Verticality Mean [deg]",
    "Verticality Median [deg]",
    "Verticality Q05 [deg]",
    "Verticality Q95 [deg]",
    "Normal Std Angle [deg]",
    "Radius [m]",
    "k-NN",
    "Sampled Points",
    "File",
    "Folder"
)

$headerRow = 1
$lastCol = $headers.Count   # 48 -> column AV
$oldLastCol = 36            # previous extent -> column AJ

# Make sure the newly-added columns (AK:AV) inherit the same bold / centered /
# bordered header formatting as the existing header cells before we write the
# new labels into them, by copying the format of the last pre-existing header
# cell (AJ1) across the extended range.
$fmtSource = $ws.Cells.Item($headerRow, $oldLastCol)
$fmtTarget = $ws.Range($ws.Cells.Item($headerRow, $oldLastCol + 1), $ws.Cells.Item($headerRow, $lastCol))
$fmtSource.Copy()
$fmtTarget.PasteSpecial(-4122)

# Write every header label (this both renames the existing columns whose text
# changed and populates the newly inserted ones).
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item($headerRow, $col).Value = $headers[$i]
}

# Match the saved selection/active cell recorded in the refreshed workbook.
[void]$ws.Range("F11").Select()
